# Add the Normal Forms analysis text to slide 6's content placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)

# --- Resize / reposition the placeholder (matches exact target EMU values) ---
$shape.Left = 42.392520904541016
$shape.Top = 96.22441101074219
$shape.Width = 812.8598022460938
$shape.Height = 395.7756042480469

# --- Build the paragraph text ---
$lines = @(
  'employee(EMPID,TicketID,fname,lname)',
  '2NF',
  'Candidate keys: EMPID, TicketID',
  'Functional dependency (EMPID->fname,lname)',
  'Functional dependency (TicketID->EMPID)',
  'device(Make,Model,Damaged,Serial_Number, Ticket_ID)',
  'BCNF',
  'Candidate keys: Serial_Number',
  'Functional dependency device(Serial_Number ->Make, Model, Damaged, Ticket_ID)',
  '',
  'customer(Fname,Lname,CellNum,Email,NotifyFreq,ReasonForVisit,TicketID,Preference) ',
  'BCNF',
  'Candidate keys: CellNum ',
  'Functional dependency (CellNum ->Fname, Lname, Email, NotifyFreq, ReasonForVisit, TicketID, Preference)'
)
$cr = [char]13

$tf = $shape.TextFrame
$tr = $tf.TextRange

# Build the text paragraph-by-paragraph via InsertAfter: this keeps each
# generated run's lang="en-US" attribute (a plain bulk `.Text =` across
# multiple paragraphs drops it on the runtime we're scripting against).
$tr.Text = $lines[0]
for ($i = 1; $i -lt $lines.Count; $i++) {
    $tr.InsertAfter($cr + $lines[$i]) | Out-Null
}

# --- Auto-shrink text to fit the placeholder ---
$tf.AutoSize = 2

# --- Set outline/indent level 1 (XML lvl="1") on the sub-bullet paragraphs ---
$lvl1Paragraphs = @(2,3,4,5,7,8,9,12,13,14)
foreach ($idx in $lvl1Paragraphs) {
    $tr.Paragraphs($idx,1).IndentLevel = 2
}

# --- Underline the candidate-key attribute names within the schema lines ---
$tr.Characters(10,5).Font.Underline = -1
$tr.Characters(183,14).Font.Underline = -1
$tr.Characters(344,7).Font.Underline = -1
